$d = $word.ActiveDocument

# 1. The lone paragraph in the body was only carrying the bookmark for
#    "_GoBack" but had accidentally been left tagged with the "Note"
#    paragraph style. Strip that so the paragraph goes back to using the
#    document's default (Normal) style, i.e. drop <w:pPr><w:pStyle .../>.
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Note") {
        $p.Style = $d.Styles.Item("Normal")
    }
}

# 2. Fix the misspelled custom style name "MarginNoteRIght" -> "MarginNoteRight"
#    (a stray capital I) so it matches the other MarginNote* styles and the
#    i18n stylesheet stops trying to load it twice.
foreach ($s in $d.Styles) {
    if ($s.NameLocal -eq "MarginNoteRIght") {
        $s.NameLocal = "MarginNoteRight"
    }
}
